# Scheduled-runner update: refresh market-board derived columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) across several sheets.
$wb = $excel.ActiveWorkbook

# ---- ALC sheet ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 4030.7778
$ws.Range("I116").Value = 1796.7273
$ws.Range("J116").Value = 7541.4287
$ws.Range("K116").Value = 1796.7273
$ws.Range("L116").Value = 7541.4287
$ws.Range("M116").Value = 1645.2727
$ws.Range("N116").Value = -14425.4287

$ws.Range("H132").Value = 2383.244
$ws.Range("I132").Value = 1405.7778
$ws.Range("J132").Value = 9421
$ws.Range("K132").Value = 4217.3334
$ws.Range("L132").Value = 28263
$ws.Range("M132").Value = -1687.3334
$ws.Range("N132").Value = -33323

$ws.Range("H137").Value = 1947.4706
$ws.Range("J137").Value = 2678.6667
$ws.Range("L137").Value = 8036.000100000001
$ws.Range("N137").Value = -13136.0001

$ws.Range("H138").Value = 2474.7827
$ws.Range("I138").Value = 2762.7273
$ws.Range("J138").Value = 2420.1724
$ws.Range("K138").Value = 8288.1819
$ws.Range("L138").Value = 7260.5172
$ws.Range("M138").Value = -3148.1819
$ws.Range("N138").Value = -17540.5172

# ---- ARM sheet ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H105").Value = 22400
$ws.Range("J105").Value = 22400
$ws.Range("L105").Value = 22400
$ws.Range("N105").Value = -29388

# ---- BSM sheet ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()

# ---- CRP sheet ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1470.75
$ws.Range("I16").Value = 1615.8572
$ws.Range("J16").Value = 1357.8889
$ws.Range("K16").Value = 1615.8572
$ws.Range("L16").Value = 1357.8889
$ws.Range("M16").Value = -1328.8572
$ws.Range("N16").Value = -1931.8889

$ws.Range("H31").Value = 202664.66
$ws.Range("I31").Value = 1523.375
$ws.Range("J31").Value = 1007229.8
$ws.Range("K31").Value = 1523.375
$ws.Range("L31").Value = 1007229.8
$ws.Range("M31").Value = -1228.375
$ws.Range("N31").Value = -1007819.8

$ws.Range("H34").Value = 202664.66
$ws.Range("I34").Value = 1523.375
$ws.Range("J34").Value = 1007229.8
$ws.Range("K34").Value = 1523.375
$ws.Range("L34").Value = 1007229.8
$ws.Range("M34").Value = -1321.375
$ws.Range("N34").Value = -1007633.8

$ws.Range("H113").Value = 1470.75
$ws.Range("I113").Value = 1615.8572
$ws.Range("J113").Value = 1357.8889
$ws.Range("K113").Value = 1615.8572
$ws.Range("L113").Value = 1357.8889
$ws.Range("M113").Value = 554.1428000000001
$ws.Range("N113").Value = -5697.8889

# ---- CUL sheet ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 504.26315
$ws.Range("I33").Value = 430
$ws.Range("J33").Value = 558.2727
$ws.Range("K33").Value = 2580
$ws.Range("L33").Value = 3349.6362
$ws.Range("M33").Value = -2297
$ws.Range("N33").Value = -3915.6362

$ws.Range("H122").Value = 3489.0344
$ws.Range("J122").Value = 3650.4363
$ws.Range("L122").Value = 32853.9267
$ws.Range("N122").Value = -37753.9267

# ---- GSM sheet ----
# Rows 125-141: market board had no listings this pass, so the derived
# price/profit columns (H:N) are cleared back out, leaving only the
# leve metadata columns (A:G).
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H125:N141").ClearContents()
